$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.148.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.233.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.12%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.67%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.04%  "

$ws.Range("E11").Value = "  -3.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.70%  "

$ws.Range("E13").Value = "  -10.05%  "

$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.869"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -11.90%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.573.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.233.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.077.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0970"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.60%  "

$ws.Range("E24").Value = "  -13.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "238.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.13%  "

$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("E29").Value = "  -1.83%  "

$ws.Range("E30").Value = "  -10.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -16.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0878"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("E39").Value = "  -7.82%  "

$ws.Range("E40").Value = "  -5.87%  "

$ws.Range("E41").Value = "  -11.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.40%  "

$ws.Range("E43").Value = "  -8.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.736.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.42%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.60%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.206"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.93%  "

$ws.Range("E49").Value = "  -10.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -16.58%  "
